# Added New Mac-Address and Document Types
# Appends a new machine record ("Machine 32") as row 33 of the
# master-machine_master sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33

$ws.Cells.Item($row, 1).Value = 10032
$ws.Cells.Item($row, 2).Value = "Machine 32"
$ws.Cells.Item($row, 3).Value = "F4-30-B9-D4-CD-6F"
$ws.Cells.Item($row, 4).Value = "FB5962911665"
$ws.Cells.Item($row, 5).Value = "192.168.0.358"
$ws.Cells.Item($row, 6).Value = 1001
$ws.Cells.Item($row, 7).Value = "eng"
$ws.Cells.Item($row, 8).Value = $true
$ws.Cells.Item($row, 9).Value = "superadmin"
$ws.Cells.Item($row, 10).Value = "now()"

# Match the final cell selection left behind by the editor in the saved file.
[void]$ws.Range("J29").Select()
